$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the columns that are being removed (E:K) - this also clears the
# custom column-width formatting that had been applied to column K.
$ws.Range("E1:K3").EntireColumn.Delete()

# Update header row
$ws.Range("A1").Value = "course_id"
$ws.Range("B1").Value = "title"
$ws.Range("C1").Value = "credits"
$ws.Range("D1").Value = "dept_name"

# Row 2
$ws.Range("A2").Value = "CCCC120001"
$ws.Range("B2").Value = "我不做人了"
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = "软件学院"

# Row 3
$ws.Range("A3").Value = "DDDD111111"
$ws.Range("B3").Value = "我死了"
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = "计算机学院"

# Reset the selection similar to the target state
$ws.Range("F4").Select()
